# The table on the slide still used the original (default) table style
# {E093A08F-DB7B-402E-965D-DEE319BD7F9F}; re-style it to
# {1AC2C49F-577F-454D-8D52-B92071D5E13C}, matching the table style picked
# from the Table Design gallery.

$p = $ppt.ActivePresentation

$targetStyleId = "{1AC2C49F-577F-454D-8D52-B92071D5E13C}"
$updated = 0

foreach ($s in $p.Slides) {
    foreach ($shp in $s.Shapes) {
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -ne $targetStyleId) {
                $tbl.ApplyStyle($targetStyleId)
                $updated = $updated + 1
            }
        }
    }
}

Write-Host "Tables restyled: $updated"
